# Refresh the forecast data: the whole table shifted one week forward and
# the forecast numbers were regenerated; the Summary sheet stats were
# recomputed to match the refreshed forecast.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Force column B (Week_Start_Date) to stay plain text so the date-like
# strings ("2025-02-02", ...) aren't auto-converted into date serials.
$wsForecast.Range("B2:B17").NumberFormat = "@"

# New weekly data: Row, Week_Start_Date, MyForecast, Amazon Mean, P70, P80, P90
$rows = @(
    @{ Row=2;  Date="2025-02-02"; D=281; E=216; F=251; G=280; H=325 },
    @{ Row=3;  Date="2025-02-09"; D=153; E=166; F=196; G=225; H=269 },
    @{ Row=4;  Date="2025-02-16"; D=170; E=168; F=199; G=229; H=274 },
    @{ Row=5;  Date="2025-02-23"; D=206; E=170; F=202; G=232; H=279 },
    @{ Row=6;  Date="2025-03-02"; D=224; E=172; F=205; G=239; H=290 },
    @{ Row=7;  Date="2025-03-09"; D=220; E=169; F=201; G=233; H=283 },
    @{ Row=8;  Date="2025-03-16"; D=216; E=166; F=200; G=235; H=291 },
    @{ Row=9;  Date="2025-03-23"; D=217; E=178; F=215; G=255; H=319 },
    @{ Row=10; Date="2025-03-30"; D=221; E=170; F=203; G=237; H=290 },
    @{ Row=11; Date="2025-04-06"; D=213; E=164; F=198; G=235; H=294 },
    @{ Row=12; Date="2025-04-13"; D=216; E=166; F=201; G=240; H=301 },
    @{ Row=13; Date="2025-04-20"; D=214; E=165; F=200; G=240; H=304 },
    @{ Row=14; Date="2025-04-27"; D=211; E=162; F=196; G=235; H=295 },
    @{ Row=15; Date="2025-05-04"; D=202; E=155; F=188; G=226; H=287 },
    @{ Row=16; Date="2025-05-11"; D=203; E=156; F=189; G=228; H=289 },
    @{ Row=17; Date="2025-05-18"; D=198; E=152; F=184; G=221; H=279 }
)

foreach ($r in $rows) {
    $wsForecast.Cells.Item($r.Row, 2).Value = $r.Date
    $wsForecast.Cells.Item($r.Row, 4).Value = $r.D
    $wsForecast.Cells.Item($r.Row, 5).Value = $r.E
    $wsForecast.Cells.Item($r.Row, 6).Value = $r.F
    $wsForecast.Cells.Item($r.Row, 7).Value = $r.G
    $wsForecast.Cells.Item($r.Row, 8).Value = $r.H
}

# Summary sheet recomputed metrics. The numeric-looking metrics are stored
# as text in this sheet, so pin the number format to text before writing
# the new values to avoid Excel re-typing them as numbers.
$textCells = @("B4", "B5", "B6", "B9", "B10", "B11", "B12", "B14")
foreach ($addr in $textCells) {
    $wsSummary.Range($addr).NumberFormat = "@"
}

$wsSummary.Range("B2").Value  = "2022-12-25 to 2025-01-26"
$wsSummary.Range("B4").Value  = "876"
$wsSummary.Range("B5").Value  = "287"
$wsSummary.Range("B6").Value  = "214"
$wsSummary.Range("B8").Value  = "32750 units"
$wsSummary.Range("B9").Value  = "3363"
$wsSummary.Range("B10").Value = "1686"
$wsSummary.Range("B11").Value = "810"
$wsSummary.Range("B12").Value = "281"
$wsSummary.Range("B14").Value = "153"
